# Actualización automática 2025-10-28 08:30:09
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------
$ws1.Range("K3").Value = 406.08
$ws1.Range("M3").Value = 1449.08

$ws1.Range("L4").Value = 886.88
$ws1.Range("M4").Value = 5211.45

$ws1.Range("L5").Value = 3282.05
$ws1.Range("M5").Value = 9711.83

$ws1.Range("P29").Value = 514.1799999999999

$ws1.Range("M36").Value = 9565.35

$ws1.Range("M52").Value = 3119.1

$ws1.Range("M53").Value = 78.41

$ws1.Range("K56").Value = "8 de 54"
$ws1.Range("L56").Value = "5 de 54"
$ws1.Range("M56").Value = "13 de 54"
$ws1.Range("P56").Value = "1 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------
$ws2.Range("F3").Value = 1855.16
$ws2.Range("F4").Value = 7718.32
$ws2.Range("F5").Value = 12993.88
$ws2.Range("F29").Value = 10838.35
$ws2.Range("F36").Value = 15328.44
$ws2.Range("F53").Value = 5098.57
$ws2.Range("F54").Value = 5098.57
$ws2.Range("F55").Value = 150.66
$ws2.Range("F56").Value = 150.66
$ws2.Range("F60").Value = 85368.86

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------
$ws3.Range("D8").Value = 514.1799999999999
$ws3.Range("E8").Value = 832.22488751609
$ws3.Range("F8").Value = 0.3818910676628507

$ws3.Range("D10").Value = 9358.65
$ws3.Range("E10").Value = -5477.570164656079
$ws3.Range("F10").Value = 2.411352097107966

$ws3.Range("D11").Value = 10130.56
$ws3.Range("E11").Value = 1700.440000000001
$ws3.Range("F11").Value = 0.8562725044374947

$ws3.Range("D12").Value = 48141.76
$ws3.Range("E12").Value = 4521.360000000001
$ws3.Range("F12").Value = 0.9141456108183488

$ws3.Range("D14").Value = 80119.63
$ws3.Range("E14").Value = 18896.87661190614
$ws3.Range("F14").Value = 0.8091542788317893
